$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 125; this shifts the existing rows 125:178 down to 126:179
# and Excel copies formatting (incl. the date number format on column D) from
# the surrounding rows automatically.
$ws.Rows("125:125").Insert()

# Populate the newly-inserted row with this week's (new) record for Haba.
$ws.Range("A125").Value = 9
$ws.Range("B125").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C125").Value = "Metropolitana"
$ws.Range("D125").Value = 44523
$ws.Range("E125").Value = 13
$ws.Range("F125").Value = 100112026
$ws.Range("G125").Value = "Haba"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 43
$ws.Range("K125").Value = 7000
$ws.Range("L125").Value = 8000
$ws.Range("M125").Value = 7512
$ws.Range("N125").Value = "$/saco 25 kilos"
$ws.Range("O125").Value = "Región del Maule"
$ws.Range("P125").Value = 300
$ws.Range("Q125").Value = 25
$ws.Range("R125").Value = "Hortaliza"
